$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 55, pushing existing rows 55-61 down to 56-62
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly record
$ws.Cells.Item(55, 1).Value = 11
$ws.Cells.Item(55, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(55, 3).Value = "Bíobío"
$ws.Cells.Item(55, 4).Value = 45154
$ws.Cells.Item(55, 4).NumberFormat = $ws.Cells.Item(56, 4).NumberFormat
$ws.Cells.Item(55, 5).Value = 8
$ws.Cells.Item(55, 6).Value = 100114007
$ws.Cells.Item(55, 7).Value = "Jengibre"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 50
$ws.Cells.Item(55, 11).Value = 18000
$ws.Cells.Item(55, 12).Value = 18000
$ws.Cells.Item(55, 13).Value = 18000
$ws.Cells.Item(55, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(55, 15).Value = "Perú"
$ws.Cells.Item(55, 16).Value = 1385
$ws.Cells.Item(55, 17).Value = 13
$ws.Cells.Item(55, 18).Value = "Hortaliza"
